$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the two new rows for "Reverse String" and "Reverse Vowels of a String"
$ws.Range("A42").Value = 344
$ws.Range("B42").Value = "Reverse String"
$ws.Range("C42").Value = "2 Pointers"

$ws.Range("A43").Value = 345
$ws.Range("B43").Value = "Reverse Vowels of a String"
$ws.Range("C43").Value = "2 Pointers/List<int>"

# Update the view to match the new selection/scroll position
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C43").Select()
